# Updates cryptos list values (Price / Volume(1h) columns) per the
# Mon Feb  5 03:51:03 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.593.81'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.284.15'
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'303.88"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').Value = "'95.74"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.58%  '
$ws.Range('E7').Value = '  -2.60%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('D10').Value = "'34.67"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('D11').Value = "'0.0780"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').Value = "'18.13"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = "'6.82"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').Value = '2.640.04'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '2.284.85'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '42.508.02'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D19').Value = "'12.90"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').Value = "'67.14"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').Value = "'235.75"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.65%  '
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').Value = "'24.64"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('E28').Value = '  +16.85%  '
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').Value = "'32.91"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').Value = "'17.80"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('D35').Value = "'4.47"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.80%  '
$ws.Range('E36').Value = '  -2.44%  '
$ws.Range('D37').Value = "'0.0682"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('E40').Value = '  -1.84%  '
$ws.Range('D41').Value = "'2.66"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('D42').Value = '1.991.63'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('E43').Value = '  -3.76%  '
$ws.Range('D44').Value = "'10.21"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = "'18.29"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.10%  '
$ws.Range('D46').Value = "'2.05"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').Value = "'2.89"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('D49').Value = "'53.53"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').Value = '2.505.52'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').Value = "'1.12"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.42%  '
